$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new material row: Material ID 902, Name "DTUMortar"
$ws.Range("A12").Value = 902
$ws.Range("B12").Value = "DTUMortar"

# Update the selected cell to match the saved view state
$ws.Range("A13").Select()
